# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - daily hourly spot prices update (next day's data)
$ws.Range("A2").Value = 45884

$ws.Range("B2").Value = 119.07
$ws.Range("C2").Value = 105.3
$ws.Range("D2").Value = 99.90000000000001
$ws.Range("E2").Value = 96.22
$ws.Range("F2").Value = 93.66
$ws.Range("G2").Value = 91.16
$ws.Range("H2").Value = 98.5
$ws.Range("I2").Value = 97.52
$ws.Range("J2").Value = 71.89
$ws.Range("K2").Value = 43.97
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 0.08
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1.72
$ws.Range("R2").Value = 23.53
$ws.Range("S2").Value = 33.23
$ws.Range("T2").Value = 61.12
$ws.Range("U2").Value = 100.82
$ws.Range("V2").Value = 111.39
$ws.Range("W2").Value = 144.75
$ws.Range("X2").Value = 142
$ws.Range("Y2").Value = 117.63
$ws.Range("Z2").Value = 69.27

# AA2 (Slot_4h_max) stays "20h-24h" - unchanged

$ws.Range("AB2").Value = 128.94
$ws.Range("AC2").Value = "22h-24h"
$ws.Range("AD2").Value = 129.82
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 128.07
$ws.Range("AG2").Value = "9h-18h"
